$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.414393096871265
$ws.Range("C2").Value = 0.07513119594463546
$ws.Range("D2").Value = 0.01161827540880012
$ws.Range("E2").Value = 0.06358328784665535
$ws.Range("F2").Value = 3.767704608165204
$ws.Range("I2").Value = 2.378409290992231
$ws.Range("J2").Value = 0.1458084058797553
$ws.Range("K2").Value = 1.245497301426099
$ws.Range("L2").Value = 0.3418581033705408
$ws.Range("M2").Value = 0.3636948708781311
$ws.Range("N2").Value = 4.004334293927457
$ws.Range("B3").Value = 1.385652285263262
$ws.Range("C3").Value = 0.06958702017958274
$ws.Range("D3").Value = 0.01145343688959599
$ws.Range("E3").Value = 0.06374172123281863
$ws.Range("F3").Value = 3.769346261606415
$ws.Range("I3").Value = 2.382981102763637
$ws.Range("J3").Value = 0.1463132986269713
$ws.Range("K3").Value = 1.212520268991398
$ws.Range("L3").Value = 0.3407908606220573
$ws.Range("M3").Value = 0.3590000546155672
$ws.Range("N3").Value = 4.024153980476534
$ws.Range("B4").Value = 1.368732865380878
$ws.Range("C4").Value = 0.06622488733900411
$ws.Range("D4").Value = 0.01135633486471832
$ws.Range("E4").Value = 0.06385006569860163
$ws.Range("F4").Value = 3.771710713770659
$ws.Range("I4").Value = 2.386574929847733
$ws.Range("J4").Value = 0.1466427731050377
$ws.Range("K4").Value = 1.192950057887913
$ws.Range("L4").Value = 0.3402659181188312
$ws.Range("M4").Value = 0.3562897929574511
$ws.Range("N4").Value = 4.037206725149581
$ws.Range("B5").Value = 1.362021437494121
$ws.Range("C5").Value = 0.06486530894528642
$ws.Range("D5").Value = 0.01131780850607278
$ws.Range("E5").Value = 0.06389700823585986
$ws.Range("F5").Value = 3.773015631818609
$ws.Range("I5").Value = 2.388237422814505
$ws.Range("J5").Value = 0.1467819426259913
$ws.Range("K5").Value = 1.185145767169672
$ws.Range("L5").Value = 0.3400848518811301
$ws.Range("M5").Value = 0.3552287840750701
$ws.Range("N5").Value = 4.042748035554332
$ws.Range("B6").Value = 1.360918097843694
$ws.Range("C6").Value = 0.06464018590206422
$ws.Range("D6").Value = 0.01131147455965298
$ws.Range("E6").Value = 0.06390497185235322
$ws.Range("F6").Value = 3.773252939854189
$ws.Range("I6").Value = 2.388525440517675
$ws.Range("J6").Value = 0.1468053482398153
$ws.Range("K6").Value = 1.183860190781644
$ws.Range("L6").Value = 0.3400567723661325
$ws.Range("M6").Value = 0.3550552312676025
$ws.Range("N6").Value = 4.043681589290287
$ws.Range("B7").Value = 1.368641609684488
$ws.Range("C7").Value = 0.06620650906393166
$ws.Range("D7").Value = 0.011355811045755
$ws.Range("E7").Value = 0.0638506874685465
$ws.Range("F7").Value = 3.771726929723016
$ws.Range("I7").Value = 2.386596549001709
$ws.Range("J7").Value = 0.1466446301156275
$ws.Range("K7").Value = 1.192844114838067
$ws.Range("L7").Value = 0.3402633430795063
$ws.Range("M7").Value = 0.3562753077818179
$ws.Range("N7").Value = 4.037280557383752
$ws.Range("B8").Value = 1.404332501681552
$ws.Range("C8").Value = 0.07321081431442167
$ws.Range("D8").Value = 0.01156059212550531
$ws.Range("E8").Value = 0.06363562419008417
$ws.Range("F8").Value = 3.767989280403299
$ws.Range("I8").Value = 2.379822457075647
$ws.Range("J8").Value = 0.1459784586040103
$ws.Range("K8").Value = 1.233986296130496
$ws.Range("L8").Value = 0.3414631078908954
$ws.Range("M8").Value = 0.3620403792188576
$ws.Range("N8").Value = 4.010984824322307
$ws.Range("B9").Value = 1.480081863560201
$ws.Range("C9").Value = 0.08728246959225316
$ws.Range("D9").Value = 0.01199436428522915
$ws.Range("E9").Value = 0.06330131184923449
$ws.Range("F9").Value = 3.771413321634995
$ws.Range("I9").Value = 2.372775455262989
$ws.Range("J9").Value = 0.1448260928000575
$ws.Range("K9").Value = 1.320037491808392
$ws.Range("L9").Value = 0.3448472720326734
$ws.Range("M9").Value = 0.3747100147999234
$ws.Range("N9").Value = 3.966424602519893
$ws.Range("B10").Value = 1.539236436787775
$ws.Range("C10").Value = 0.0978309879610606
$ws.Range("D10").Value = 0.0123321862371597
$ws.Range("E10").Value = 0.06310849629421522
$ws.Range("F10").Value = 3.78047682114493
$ws.Range("I10").Value = 2.371395330883885
$ws.Range("J10").Value = 0.144072664776548
$ws.Range("K10").Value = 1.386533597482497
$ws.Range("L10").Value = 0.3479592734738333
$ws.Range("M10").Value = 0.3848471188867819
$ws.Range("N10").Value = 3.937952070596324
$ws.Range("B11").Value = 1.566906207619013
$ws.Range("C11").Value = 0.1026765632283002
$ws.Range("D11").Value = 0.01248992391124659
$ws.Range("E11").Value = 0.06303214242176036
$ws.Range("F11").Value = 3.786020607849508
$ws.Range("I11").Value = 2.371591185374314
$ws.Range("J11").Value = 0.143750011023883
$ws.Range("K11").Value = 1.417495914043229
$ws.Range("L11").Value = 0.3495102468375251
$ws.Range("M11").Value = 0.3896381665626407
$ws.Range("N11").Value = 3.925924507755767
$ws.Range("B12").Value = 1.577493034188706
$ws.Range("C12").Value = 0.1045182850903075
$ws.Range("D12").Value = 0.01255022999209388
$ws.Range("E12").Value = 0.06300485413342471
$ws.Range("F12").Value = 3.788324018584603
$ws.Range("I12").Value = 2.371783695792701
$ws.Range("J12").Value = 0.143630707940881
$ws.Range("K12").Value = 1.429322913603613
$ws.Range("L12").Value = 0.350116954193382
$ws.Range("M12").Value = 0.3914781673295025
$ws.Range("N12").Value = 3.921502910547275
$ws.Range("B13").Value = 1.575208135486548
$ws.Range("C13").Value = 0.1041213337799149
$ws.Range("D13").Value = 0.01253721661163354
$ws.Range("E13").Value = 0.06301065899922076
$ws.Range("F13").Value = 3.787818863707514
$ws.Range("I13").Value = 2.371736973379143
$ws.Range("J13").Value = 0.1436562741010636
$ws.Range("K13").Value = 1.42677121469481
$ws.Range("L13").Value = 0.3499854277551435
$ws.Range("M13").Value = 0.3910807471823361
$ws.Range("N13").Value = 3.922449266193539
$ws.Range("B14").Value = 1.567775012338871
$ws.Range("C14").Value = 0.1028279460768999
$ws.Range("D14").Value = 0.01249487387466885
$ws.Range("E14").Value = 0.06302986487150974
$ws.Range("F14").Value = 3.786206021541332
$ws.Range("I14").Value = 2.371604652082326
$ws.Range("J14").Value = 0.1437401382427126
$ws.Range("K14").Value = 1.418466880879237
$ws.Range("L14").Value = 0.3495597728626194
$ws.Range("M14").Value = 0.3897890293091351
$ws.Range("N14").Value = 3.925558075141879
$ws.Range("B15").Value = 1.563236177671342
$ws.Range("C15").Value = 0.1020365967794135
$ws.Range("D15").Value = 0.01246901222572561
$ws.Range("E15").Value = 0.0630418404401567
$ws.Range("F15").Value = 3.785244682885221
$ws.Range("I15").Value = 2.371539010511825
$ws.Range("J15").Value = 0.1437918820951154
$ws.Range("K15").Value = 1.413393543551138
$ws.Range("L15").Value = 0.3493015697001596
$ws.Range("M15").Value = 0.3890011634523987
$ws.Range("N15").Value = 3.927479628473151
$ws.Range("B16").Value = 1.537443411343048
$ws.Range("C16").Value = 0.09751526818618572
$ws.Range("D16").Value = 0.0123219585895491
$ws.Range("E16").Value = 0.06311371407303046
$ws.Range("F16").Value = 3.780143093348471
$ws.Range("I16").Value = 2.371399097304412
$ws.Range("J16").Value = 0.1440941543174361
$ws.Range("K16").Value = 1.38452445937935
$ws.Range("L16").Value = 0.3478606304908993
$ws.Range("M16").Value = 0.3845376190677783
$ws.Range("N16").Value = 3.938756708394195
$ws.Range("B17").Value = 1.521814760773111
$ws.Range("C17").Value = 0.09475365245425849
$ws.Range("D17").Value = 0.01223277906380105
$ws.Range("E17").Value = 0.06316070969875476
$ws.Range("F17").Value = 3.777377176514619
$ws.Range("I17").Value = 2.371524157595402
$ws.Range("J17").Value = 0.1442847260495235
$ws.Range("K17").Value = 1.366996608479354
$ws.Range("L17").Value = 0.347011266862097
$ws.Range("M17").Value = 0.3818453182075032
$ws.Range("N17").Value = 3.945911680021368
$ws.Range("B18").Value = 1.512897139010875
$ws.Range("C18").Value = 0.09316966517822323
$ws.Range("D18").Value = 0.01218186814721633
$ws.Range("E18").Value = 0.06318880997465204
$ws.Range("F18").Value = 3.775920034399419
$ws.Range("I18").Value = 2.371673623390066
$ws.Range("J18").Value = 0.1443962288574614
$ws.Range("K18").Value = 1.356982170309095
$ws.Range("L18").Value = 0.346535475365485
$ws.Range("M18").Value = 0.3803136910709881
$ws.Range("N18").Value = 3.950114070315564
$ws.Range("B19").Value = 1.509890089377052
$ws.Range("C19").Value = 0.09263411213575523
$ws.Range("D19").Value = 0.01216469662791297
$ws.Range("E19").Value = 0.06319850820475192
$ws.Range("F19").Value = 3.77544964740261
$ws.Range("I19").Value = 2.371737549639242
$ws.Range("J19").Value = 0.1444343068544161
$ws.Range("K19").Value = 1.353602989561466
$ws.Range("L19").Value = 0.3463765707806914
$ws.Range("M19").Value = 0.3797980161985564
$ws.Range("N19").Value = 3.951551876601471
$ws.Range("B20").Value = 1.523471053613235
$ws.Range("C20").Value = 0.09504717324706746
$ws.Range("D20").Value = 0.01224223283035997
$ws.Range("E20").Value = 0.0631555962901782
$ws.Range("F20").Value = 3.777657772966151
$ws.Range("I20").Value = 2.371502820914785
$ws.Range("J20").Value = 0.1442642437202544
$ws.Range("K20").Value = 1.368855533644734
$ws.Range("L20").Value = 0.3471003650358568
$ws.Range("M20").Value = 0.3821301684794634
$ws.Range("N20").Value = 3.945141012748081
$ws.Range("B21").Value = 1.569955348942699
$ws.Range("C21").Value = 0.1032076602479606
$ws.Range("D21").Value = 0.01250729545467166
$ws.Range("E21").Value = 0.06302417959642526
$ws.Range("F21").Value = 3.786674214794971
$ws.Range("I21").Value = 2.371640306997321
$ws.Range("J21").Value = 0.1437154272524879
$ws.Range("K21").Value = 1.42090329235441
$ws.Range("L21").Value = 0.3496842724075435
$ws.Range("M21").Value = 0.3901677404762012
$ws.Range("N21").Value = 3.924641333507438
$ws.Range("B22").Value = 1.600969999031065
$ws.Range("C22").Value = 0.1085807134527386
$ws.Range("D22").Value = 0.0126838719714506
$ws.Range("E22").Value = 0.06294776075778863
$ws.Range("F22").Value = 3.793756534590628
$ws.Range("I22").Value = 2.372419920501258
$ws.Range("J22").Value = 0.1433735201175157
$ws.Range("K22").Value = 1.455515332997862
$ws.Range("L22").Value = 0.3514859778502313
$ws.Range("M22").Value = 0.3955707178915162
$ws.Range("N22").Value = 3.912018721001644
$ws.Range("B23").Value = 1.584358976894691
$ws.Range("C23").Value = 0.1057093647563363
$ws.Range("D23").Value = 0.0125893270558386
$ws.Range("E23").Value = 0.06298768313292857
$ws.Range("F23").Value = 3.789867789205204
$ws.Range("I23").Value = 2.37194074539822
$ws.Range("J23").Value = 0.1435544703982856
$ws.Range("K23").Value = 1.436987807138365
$ws.Range("L23").Value = 0.3505140596673044
$ws.Range("M23").Value = 0.3926733571748926
$ws.Range("N23").Value = 3.918684721567999
$ws.Range("B24").Value = 1.522722033338255
$ws.Range("C24").Value = 0.09491446096875222
$ws.Range("D24").Value = 0.01223795766197711
$ws.Range("E24").Value = 0.06315790469205851
$ws.Range("F24").Value = 3.777530500928677
$ws.Range("I24").Value = 2.371512225592213
$ws.Range("J24").Value = 0.1442734977359503
$ws.Range("K24").Value = 1.368014918697611
$ws.Range("L24").Value = 0.3470600447513306
$ws.Range("M24").Value = 0.3820013371638069
$ws.Range("N24").Value = 3.94548915447416
$ws.Range("B25").Value = 1.458974061452636
$ws.Range("C25").Value = 0.08343919201762162
$ws.Range("D25").Value = 0.01187362266945158
$ws.Range("E25").Value = 0.06338244389136793
$ws.Range("F25").Value = 3.769336608620051
$ws.Range("I25").Value = 2.374014663730335
$ws.Range("J25").Value = 0.1451214210619156
$ws.Range("K25").Value = 1.29618335090791
$ws.Range("L25").Value = 0.343821669937725
$ws.Range("M25").Value = 0.3711367947434852
$ws.Range("N25").Value = 3.977729850541863
